$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add three new names to column A, below the existing list.
# The entry order reproduces the order new unique strings were recorded
# in the shared-strings table: A14, then A15, then A13.
$ws.Range("A14").Value = "Justine Faith"
$ws.Range("A15").Value = "Timothy Jonah"
$ws.Range("A13").Value = "Justin John"

# Reflect the final active cell/selection left after the edits.
$ws.Range("G11").Select()
